# "encrypt link and change register users"
# Update the registered-user rows (2-11) on Hoja1:
#   - Grupos (B): Matecero -> Ciencias
#   - Nombres (C): append a surname/suffix to each name
#   - Usuario (D): replace numeric ids with ab10..ab100 text codes (right aligned like column H)
#   - Ciclo (H): ANUAL -> VERANO
# Also move the sheet selection to B15, matching the saved UI state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "Aquilino robert",
    "Eloy ticona",
    "Pascual  perez",
    "Pedro Miguel  lopez",
    "Cristian Acev tipo",
    "Francisco Luis  juarez",
    "Marti Cha  perez",
    "Luis Javier Pos lipo",
    "Ivan Sevil  tica",
    "Damian Mon  morales"
)

$usuarios = @("ab10","ab20","ab30","ab40","ab50","ab60","ab70","ab80","ab90","ab100")

for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 2

    $ws.Cells.Item($row, 2).Value = "Ciencias"
    $ws.Cells.Item($row, 3).Value = $names[$i]

    $ws.Cells.Item($row, 4).Value = $usuarios[$i]
    $ws.Cells.Item($row, 4).HorizontalAlignment = -4152

    $ws.Cells.Item($row, 8).Value = "VERANO"
}

$ws.Range("B15").Select()
